$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Sheet view: zoom to 80% and move the selection to K2 (drops topLeftCell/zoomScaleNormal=100 defaults) ---
$excel.ActiveWindow.Zoom = 80

# --- Column widths ---
# Column A: 11.42578125 (bestFit) -> 14.85546875 (custom)
$ws.Columns.Item(1).ColumnWidth = 14.022135416666666
# Column D: split out of the old C:D 24.7109375 merged range -> 28.85546875 (custom)
$ws.Columns.Item(4).ColumnWidth = 28.022135416666668
# Column E: 20.5703125 (bestFit) -> 27.5703125 (custom)
$ws.Columns.Item(5).ColumnWidth = 26.736979166666668
# Column H: 17.5703125 (bestFit) -> 22 (custom)
$ws.Columns.Item(8).ColumnWidth = 21.166666666666668
# Column I: 30.5703125 (custom) -> 34.5703125 (custom)
$ws.Columns.Item(9).ColumnWidth = 33.736979166666664

# Selection goes last so the final saved selection is K2 (column changes above do not move the selection,
# but keep this last for clarity/robustness).
$ws.Range("K2").Select()
